{"js": "// Replace the two-digit multiplication problems in the document body.\n// Each key appears exactly once in the document, so a straightforward\n// search + replace of the exact \"NN\u00d7NN=\" text is safe and unambiguous.\nconst replacements = [\n  [\"27\u00d716=\", \"75\u00d774=\"],\n  [\"31\u00d767=\", \"98\u00d742=\"],\n  [\"80\u00d739=\", \"95\u00d723=\"],\n  [\"16\u00d784=\", \"76\u00d772=\"],\n  [\"38\u00d755=\", \"36\u00d784=\"],\n  [\"48\u00d798=\", \"90\u00d756=\"],\n  [\"75\u00d756=\", \"15\u00d718=\"],\n  [\"30\u00d721=\", \"31\u00d755=\"],\n  [\"67\u00d798=\", \"37\u00d768=\"],\n  [\"30\u00d748=\", \"17\u00d769=\"],\n  [\"93\u00d741=\", \"27\u00d758=\"],\n  [\"71\u00d764=\", \"30\u00d726=\"],\n  [\"95\u00d774=\", \"34\u00d797=\"],\n  [\"79\u00d712=\", \"32\u00d768=\"],\n  [\"18\u00d728=\", \"18\u00d752=\"],\n  [\"38\u00d725=\", \"63\u00d751=\"],\n  [\"47\u00d745=\", \"18\u00d766=\"],\n  [\"75\u00d779=\", \"40\u00d784=\"],\n  [\"79\u00d767=\", \"78\u00d743=\"],\n  [\"99\u00d713=\", \"73\u00d715=\"],\n  [\"11\u00d781=\", \"14\u00d757=\"],\n  [\"54\u00d797=\", \"75\u00d742=\"],\n  [\"65\u00d742=\", \"94\u00d722=\"],\n  [\"26\u00d750=\", \"56\u00d769=\"],\n  [\"68\u00d721=\", \"98\u00d721=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems throughout the document.\n# Each \"NN\u00d7NN=\" string occurs exactly once, so a plain Find/Replace\n# (whole document scope, no wildcards needed) is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"27\u00d716=\", \"75\u00d774=\"),\n    @(\"31\u00d767=\", \"98\u00d742=\"),\n    @(\"80\u00d739=\", \"95\u00d723=\"),\n    @(\"16\u00d784=\", \"76\u00d772=\"),\n    @(\"38\u00d755=\", \"36\u00d784=\"),\n    @(\"48\u00d798=\", \"90\u00d756=\"),\n    @(\"75\u00d756=\", \"15\u00d718=\"),\n    @(\"30\u00d721=\", \"31\u00d755=\"),\n    @(\"67\u00d798=\", \"37\u00d768=\"),\n    @(\"30\u00d748=\", \"17\u00d769=\"),\n    @(\"93\u00d741=\", \"27\u00d758=\"),\n    @(\"71\u00d764=\", \"30\u00d726=\"),\n    @(\"95\u00d774=\", \"34\u00d797=\"),\n    @(\"79\u00d712=\", \"32\u00d768=\"),\n    @(\"18\u00d728=\", \"18\u00d752=\"),\n    @(\"38\u00d725=\", \"63\u00d751=\"),\n    @(\"47\u00d745=\", \"18\u00d766=\"),\n    @(\"75\u00d779=\", \"40\u00d784=\"),\n    @(\"79\u00d767=\", \"78\u00d743=\"),\n    @(\"99\u00d713=\", \"73\u00d715=\"),\n    @(\"11\u00d781=\", \"14\u00d757=\"),\n    @(\"54\u00d797=\", \"75\u00d742=\"),\n    @(\"65\u00d742=\", \"94\u00d722=\"),\n    @(\"26\u00d750=\", \"56\u00d769=\"),\n    @(\"68\u00d721=\", \"98\u00d721=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
